# Geopackages.xlsx -- citynames / area-name sheets, style tweaks
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- new sheets -----------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

# ===========================================================================
# Sheet2 - "Plaatsnamen" (city names)
# ===========================================================================
$ws2.Range("A1").Value = "Plaatsnamen"

$ws2.Range("A4").Value = "Schaal"
$ws2.Range("B4").Value = "Layer"
$ws2.Range("C4").Value = "punt"
$ws2.Range("D4").Value = "vlak"
$ws2.Range("E4").Value = "multivlak"
$ws2.Range("F4").Value = "Totaal"
$ws2.Range("H4").Value = "punt"
$ws2.Range("I4").Value = "vlak"

$ws2.Range("A5").Value = 10
$ws2.Range("C5").Value = 0
$ws2.Range("D5").Value = 2589
$ws2.Range("E5").Value = 158
$ws2.Range("F5").Formula = "=SUM(C5:E5)"

# the three labels below are entered out of row order on purpose, so the
# shared-string table receives new entries in the same sequence the
# original authoring tool produced them (>25000, >5000, <10000)
$ws2.Range("I8").Value = ">25000"
$ws2.Range("I7").Value = ">5000"
$ws2.Range("H6").Value = "<10000"

$ws2.Range("A6").Value = 250
$ws2.Range("C6").Value = 1949
$ws2.Range("D6").Value = 659
$ws2.Range("F6").Formula = "=SUM(C6:E6)"

$ws2.Range("A7").Value = 500
$ws2.Range("C7").Value = 699
$ws2.Range("D7").Value = 659
$ws2.Range("F7").Formula = "=SUM(C7:E7)"

$ws2.Range("A8").Value = 1000
$ws2.Range("C8").Value = 1376
$ws2.Range("D8").Value = 121
$ws2.Range("F8").Formula = "=SUM(C8:E8)"

$ws2.Range("I6").Select()

# ===========================================================================
# Sheet3 - "Plaat tekstgrootte" (label text size)
# ===========================================================================
$ws3.Range("A3").Value = "Plaat tekstgrootte"

$ws3.Range("C4").Value = 26
$ws3.Range("E4").Value = 10
$ws3.Range("F4").Value = 100

$ws3.Range("C5").Value = 21
$ws3.Range("E5").Formula = "=`$C5/`$C`$4*E`$4"
$ws3.Range("F5").Formula = "=`$C5/`$C`$4*F`$4"

$ws3.Range("C6").Value = 18
$ws3.Range("E6").Formula = "=`$C6/`$C`$4*E`$4"
$ws3.Range("F6").Formula = "=`$C6/`$C`$4*F`$4"

$ws3.Range("C7").Value = 14
$ws3.Range("E7").Formula = "=`$C7/`$C`$4*E`$4"
$ws3.Range("F7").Formula = "=`$C7/`$C`$4*F`$4"

$ws3.Range("C8").Value = 12
$ws3.Range("E8").Formula = "=`$C8/`$C`$4*E`$4"
$ws3.Range("F8").Formula = "=`$C8/`$C`$4*F`$4"

$ws3.Range("C9").Value = 9.5
$ws3.Range("E9").Formula = "=`$C9/`$C`$4*E`$4"
$ws3.Range("F9").Formula = "=`$C9/`$C`$4*F`$4"

# number format / style for the E2:F9 "Comma" block
$tmp = $ws3.Range("Z100")
$tmp.Style = "Comma"
$tmp.Style = "Normal"
$ws3.Range("E2:F9").NumberFormat = "_-* #,##0.0_-;\-* #,##0.0_-;_-* ""-""??_-;_-@_-"

$ws3.Range("A4").Select()
$ws3.Activate()

# --- workbook-level bits ---------------------------------------------------
$wb.Worksheets.Item(1).Select()
$ws3.Activate()

Write-Host "done"
